$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Total" column (I) for the new 2024 data.
# This shifts the old "Total" column (I) to J, carrying its formulas/refs along.
$ws.Columns("I").Insert()

# Match formatting (style) of the neighboring year columns for the new column.
$ws.Range("H1:H14").Copy()
$ws.Range("I1:I14").PasteSpecial(-4122)

# New header for 2024 column.
$ws.Range("I1").Value = 2024

# Updated raw counts for existing years (per source data correction).
$ws.Range("G2").Value = 2
$ws.Range("G4").Value = 8

# New 2024 data for each age-group row.
$ws.Range("I2").Value = 3
$ws.Range("I3").Value = 2
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 9
$ws.Range("I8").Value = 5
$ws.Range("I9").Value = 5
$ws.Range("I10").Value = 5
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 8
$ws.Range("I13").Value = 0
$ws.Range("I14").Value = 60

# Extend each row's running-total formula (old "Total" column, now J) to include the new I column.
$ws.Range("J2").Formula = "=SUM(B2:I2)"
$ws.Range("J3").Formula = "=SUM(B3:I3)"
$ws.Range("J4:J14").Formula = "=SUM(B4:I4)"

# Leave the selection where the author last left it.
$ws.Range("I16").Select()
